# The disposition code abbreviations (CONV, ACQ, DIS, ...) that used to live
# in column A have been retired -- the "new disposition codes" are simply the
# full description text, matching column B. Update column A for each data
# row (4-18) so that it mirrors column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CourtDisp")

$ws.Cells.Item(4, 1).Value = "Convicted"
$ws.Cells.Item(5, 1).Value = "Probation Without Verdict"
$ws.Cells.Item(6, 1).Value = "Not Guilty by Reason of Insanity"
$ws.Cells.Item(7, 1).Value = "Acquitted"
$ws.Cells.Item(8, 1).Value = "Dismissed"
$ws.Cells.Item(9, 1).Value = "Civil Procedure"
$ws.Cells.Item(10, 1).Value = "Off Calendar"
$ws.Cells.Item(11, 1).Value = "Guilty But Mentally Ill"
$ws.Cells.Item(12, 1).Value = "Transferred to Juvenile Court"
$ws.Cells.Item(13, 1).Value = "Mistrial"
$ws.Cells.Item(14, 1).Value = "Nolle Prosequi"
$ws.Cells.Item(15, 1).Value = "Other"
$ws.Cells.Item(16, 1).Value = "Extradited"
$ws.Cells.Item(17, 1).Value = "Not Disposition By Court"
$ws.Cells.Item(18, 1).Value = "Missing/Unknown"

# Leave the selection where it ended up when the author last saved the file.
$ws.Range("E18").Select()
